$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 3185.0833
$ws.Cells.Item(40, 9).Value = 3201.7273
$ws.Cells.Item(40, 10).Value = 3002
$ws.Cells.Item(40, 11).Value = 3201.7273
$ws.Cells.Item(40, 12).Value = 3002
$ws.Cells.Item(40, 13).Value = -3026.7273
$ws.Cells.Item(40, 14).Value = -3352
$ws.Cells.Item(137, 8).Value = 48675.95
$ws.Cells.Item(137, 9).Value = 60716.867
$ws.Cells.Item(137, 10).Value = 3522.5
$ws.Cells.Item(137, 11).Value = 182150.601
$ws.Cells.Item(137, 12).Value = 10567.5
$ws.Cells.Item(137, 13).Value = -179600.601
$ws.Cells.Item(137, 14).Value = -15667.5
$ws.Cells.Item(138, 8).Value = 3324.0781
$ws.Cells.Item(138, 10).Value = 3510.3208
$ws.Cells.Item(138, 12).Value = 10530.9624
$ws.Cells.Item(138, 14).Value = -20810.9624
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 8888.909
$ws.Cells.Item(32, 9).Value = 5022.019
$ws.Cells.Item(32, 10).Value = 24653.924
$ws.Cells.Item(32, 11).Value = 5022.019
$ws.Cells.Item(32, 12).Value = 24653.924
$ws.Cells.Item(32, 13).Value = -4735.019
$ws.Cells.Item(32, 14).Value = -25227.924
$ws.Cells.Item(61, 8).Value = 7534.875
$ws.Cells.Item(61, 9).Value = 10872.25
$ws.Cells.Item(61, 11).Value = 10872.25
$ws.Cells.Item(61, 13).Value = -10660.25
$ws.Cells.Item(119, 8).Value = 55486.4
$ws.Cells.Item(119, 10).Value = 55486.4
$ws.Cells.Item(119, 12).Value = 55486.4
$ws.Cells.Item(119, 14).Value = -65162.4
$ws.Cells.Item(122, 8).Value = 996903.8
$ws.Cells.Item(122, 9).Value = 3253.111
$ws.Cells.Item(122, 11).Value = 9759.332999999999
$ws.Cells.Item(122, 13).Value = -7309.332999999999
$ws.Cells.Item(136, 8).Value = 7534.875
$ws.Cells.Item(136, 9).Value = 10872.25
$ws.Cells.Item(136, 11).Value = 32616.75
$ws.Cells.Item(136, 13).Value = -30066.75
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 2848.6667
$ws.Cells.Item(22, 9).Value = 3119.7144
$ws.Cells.Item(22, 11).Value = 3119.7144
$ws.Cells.Item(22, 13).Value = -2946.7144
$ws.Cells.Item(134, 8).Value = 3139.8096
$ws.Cells.Item(134, 9).Value = 1565.6061
$ws.Cells.Item(134, 10).Value = 8911.888999999999
$ws.Cells.Item(134, 11).Value = 4696.8183
$ws.Cells.Item(134, 12).Value = 26735.667
$ws.Cells.Item(134, 13).Value = -2161.8183
$ws.Cells.Item(134, 14).Value = -31805.667
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(50, 8).Value = 8888.777
$ws.Cells.Item(50, 10).Value = 8888.777
$ws.Cells.Item(50, 12).Value = 8888.777
$ws.Cells.Item(50, 14).Value = -10138.777
$ws.Cells.Item(51, 8).Value = 30741.8
$ws.Cells.Item(51, 9).Value = 700
$ws.Cells.Item(51, 10).Value = 38252.25
$ws.Cells.Item(51, 11).Value = 700
$ws.Cells.Item(51, 12).Value = 38252.25
$ws.Cells.Item(51, 13).Value = 36
$ws.Cells.Item(51, 14).Value = -39724.25
$ws.Cells.Item(61, 8).Value = 30741.8
$ws.Cells.Item(61, 9).Value = 700
$ws.Cells.Item(61, 10).Value = 38252.25
$ws.Cells.Item(61, 11).Value = 700
$ws.Cells.Item(61, 12).Value = 38252.25
$ws.Cells.Item(61, 13).Value = -352
$ws.Cells.Item(61, 14).Value = -38948.25
$ws.Cells.Item(132, 8).Value = 129250.93
$ws.Cells.Item(132, 9).Value = 85989.664
$ws.Cells.Item(132, 11).Value = 257968.992
$ws.Cells.Item(132, 13).Value = -255438.992
$ws.Cells.Item(134, 8).Value = 2406.3684
$ws.Cells.Item(134, 9).Value = 1623.326
$ws.Cells.Item(134, 11).Value = 4869.978
$ws.Cells.Item(134, 13).Value = -2334.978
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(11, 8).Value = 7298.2
$ws.Cells.Item(11, 10).Value = 7403.5
$ws.Cells.Item(11, 12).Value = 22210.5
$ws.Cells.Item(11, 14).Value = -22490.5
$ws.Cells.Item(37, 8).Value = 44377.777
$ws.Cells.Item(37, 10).Value = 44377.777
$ws.Cells.Item(37, 12).Value = 133133.331
$ws.Cells.Item(37, 14).Value = -133357.331
$ws.Cells.Item(64, 8).Value = 0
$ws.Cells.Item(64, 9).Value = 0
$ws.Cells.Item(64, 11).Value = 0
$ws.Cells.Item(64, 13).ClearContents()
$ws.Cells.Item(67, 8).Value = 0
$ws.Cells.Item(67, 9).Value = 0
$ws.Cells.Item(67, 11).Value = 0
$ws.Cells.Item(67, 13).ClearContents()
$ws.Cells.Item(98, 8).Value = 1699.1428
$ws.Cells.Item(98, 10).Value = 1872.8182
$ws.Cells.Item(98, 12).Value = 5618.4546
$ws.Cells.Item(98, 14).Value = -8614.454600000001
$ws.Cells.Item(104, 8).Value = 2500
$ws.Cells.Item(104, 9).Value = 0
$ws.Cells.Item(104, 10).Value = 2500
$ws.Cells.Item(104, 11).Value = 0
$ws.Cells.Item(104, 12).Value = 7500
$ws.Cells.Item(104, 13).ClearContents()
$ws.Cells.Item(104, 14).Value = -12742
$ws.Cells.Item(107, 8).Value = 1050.8572
$ws.Cells.Item(107, 9).Value = 192.66667
$ws.Cells.Item(107, 11).Value = 578.00001
$ws.Cells.Item(107, 13).Value = 1341.99999
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(49, 8).Value = 33356666
$ws.Cells.Item(49, 10).Value = 35000
$ws.Cells.Item(49, 12).Value = 35000
$ws.Cells.Item(49, 14).Value = -35368
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(13, 8).Value = 2000
$ws.Cells.Item(13, 10).Value = 2000
$ws.Cells.Item(13, 12).Value = 2000
$ws.Cells.Item(13, 14).Value = -2280
$ws.Cells.Item(42, 8).Value = 333339680
$ws.Cells.Item(42, 9).Value = 333339680
$ws.Cells.Item(42, 11).Value = 333339680
$ws.Cells.Item(42, 13).Value = -333339117
$ws.Cells.Item(49, 8).Value = 333339680
$ws.Cells.Item(49, 9).Value = 333339680
$ws.Cells.Item(49, 11).Value = 333339680
$ws.Cells.Item(49, 13).Value = -333339533
$ws.Cells.Item(61, 8).Value = 6946321
$ws.Cells.Item(61, 9).Value = 8547943
$ws.Cells.Item(61, 11).Value = 8547943
$ws.Cells.Item(61, 13).Value = -8547741
$ws.Cells.Item(88, 8).Value = 0
$ws.Cells.Item(88, 10).Value = 0
$ws.Cells.Item(88, 12).Value = 0
$ws.Cells.Item(88, 14).ClearContents()
$ws.Cells.Item(91, 8).Value = 0
$ws.Cells.Item(91, 10).Value = 0
$ws.Cells.Item(91, 12).Value = 0
$ws.Cells.Item(91, 14).ClearContents()
$ws.Cells.Item(113, 8).Value = 6946321
$ws.Cells.Item(113, 9).Value = 8547943
$ws.Cells.Item(113, 11).Value = 8547943
$ws.Cells.Item(113, 13).Value = -8545773
$ws.Cells.Item(122, 8).Value = 7187.364
$ws.Cells.Item(122, 9).Value = 2092.3333
$ws.Cells.Item(122, 10).Value = 9098
$ws.Cells.Item(122, 11).Value = 6276.999899999999
$ws.Cells.Item(122, 12).Value = 27294
$ws.Cells.Item(122, 13).Value = -3826.999899999999
$ws.Cells.Item(122, 14).Value = -32194
$ws.Cells.Item(132, 8).Value = 9537.087
$ws.Cells.Item(132, 9).Value = 12169.667
$ws.Cells.Item(132, 10).Value = 4601
$ws.Cells.Item(132, 11).Value = 36509.001
$ws.Cells.Item(132, 12).Value = 13803
$ws.Cells.Item(132, 13).Value = -33979.001
$ws.Cells.Item(132, 14).Value = -18863
$ws.Cells.Item(136, 8).Value = 29958.525
$ws.Cells.Item(136, 9).Value = 40534.85
$ws.Cells.Item(136, 10).Value = 7992.3076
$ws.Cells.Item(136, 11).Value = 121604.55
$ws.Cells.Item(136, 12).Value = 23976.9228
$ws.Cells.Item(136, 13).Value = -119054.55
$ws.Cells.Item(136, 14).Value = -29076.9228
$ws.Cells.Item(138, 8).Value = 82123.5
$ws.Cells.Item(138, 10).Value = 82123.5
$ws.Cells.Item(138, 12).Value = 82123.5
$ws.Cells.Item(138, 14).Value = -92403.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 2419.3
$ws.Cells.Item(122, 9).Value = 2419.3
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 7257.900000000001
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -4807.900000000001
$ws.Cells.Item(122, 14).ClearContents()
$ws.Cells.Item(132, 8).Value = 40041824
$ws.Cells.Item(132, 9).Value = 71438320
$ws.Cells.Item(132, 10).Value = 82648.37
$ws.Cells.Item(132, 11).Value = 214314960
$ws.Cells.Item(132, 12).Value = 247945.11
$ws.Cells.Item(132, 13).Value = -214312430
$ws.Cells.Item(132, 14).Value = -253005.11
